$wb = $excel.ActiveWorkbook

# ============================================================
# Sheet: SCHEME_MEASURES  (MQMS01..05 -> MQME001..005)
# ============================================================
$ws = $wb.Worksheets.Item("SCHEME_MEASURES")
$ws.Range("A2").Value = "MQME001"
$ws.Range("A3").Value = "MQME002"
$ws.Range("A4").Value = "MQME003"
$ws.Range("A5").Value = "MQME004"
$ws.Range("A6").Value = "MQME005"

# ============================================================
# Sheet: METADATA_ISSUES  (rule codes renumbered)
# ============================================================
$ws = $wb.Worksheets.Item("METADATA_ISSUES")
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 1).Value = "MQME014"
}
for ($r = 9; $r -le 93; $r++) {
    $ws.Cells.Item($r, 1).Value = "MQME008"
}
for ($r = 94; $r -le 95; $r++) {
    $ws.Cells.Item($r, 1).Value = "MQME009"
}
$ws.Cells.Item(96, 1).Value = "MQME010"
for ($r = 97; $r -le 98; $r++) {
    $ws.Cells.Item($r, 1).Value = "MQME011"
}

# ============================================================
# Sheet: METADATA_MEASURES
#   Remove row "MQME00 - Total number of columns" entirely,
#   renumber remaining two rows to MQME006 / MQME007
# ============================================================
$ws = $wb.Worksheets.Item("METADATA_MEASURES")
$ws.Rows.Item(2).Delete()
$ws.Range("A2").Value = "MQME006"
$ws.Range("A3").Value = "MQME007"

# ============================================================
# Sheet: METADATA_METRICS
#   Grows from 7 data rows (IQME01..07) to 11 data rows (MQID001..011)
#   Insert 4 new rows before row 2, then rewrite all data rows 2..12.
# ============================================================
$ws = $wb.Worksheets.Item("METADATA_METRICS")
$ws.Range("A2:C5").Insert()
$ws.Range("A2:C12").ClearFormats()

$ws.Range("A2").Value = "MQID001"
$ws.Range("B2").Value = "Table names in singular"
$ws.Range("C2").Value = "'100.00%"

$ws.Range("A3").Value = "MQID002"
$ws.Range("B3").Value = "Table with recommended name length"
$ws.Range("C3").Value = "'100.00%"

$ws.Range("A4").Value = "MQID003"
$ws.Range("B4").Value = "Columns with correct prefixes"
$ws.Range("C4").Value = "'99.33%"

$ws.Range("A5").Value = "MQID004"
$ws.Range("B5").Value = "Columns with recommended name size"
$ws.Range("C5").Value = "'100.00%"

$ws.Range("A6").Value = "MQID005"
$ws.Range("B6").Value = "Columns with comments"
$ws.Range("C6").Value = "'91.81%"

$ws.Range("A7").Value = "MQID006"
$ws.Range("B7").Value = "Table with standard PK prefixes"
$ws.Range("C7").Value = "'97.83%"

$ws.Range("A8").Value = "MQID007"
$ws.Range("B8").Value = "Table with standard FK prefixes"
$ws.Range("C8").Value = "'99.48%"

$ws.Range("A9").Value = "MQID008"
$ws.Range("B9").Value = "Table with standard UK prefixes"
$ws.Range("C9").Value = "'88.89%"

$ws.Range("A10").Value = "MQID009"
$ws.Range("B10").Value = "NUMBER columns with valid scale"
$ws.Range("C10").Value = "'100.00%"

$ws.Range("A11").Value = "MQID010"
$ws.Range("B11").Value = "Columns with valid num_distinct"
$ws.Range("C11").Value = "'100.00%"

$ws.Range("A12").Value = "MQID011"
$ws.Range("B12").Value = "Columns with valid num_nulls"
$ws.Range("C12").Value = "'100.00%"
